# "fix class schedule, update performance"
#
# 1) Fix class schedule: FALL BREAK was recorded a class late. Row 14
#    (Tues, Dec 3) held the "Pandas: Reshaping" class content while row 15
#    (Thurs, Dec 5) held "FALL BREAK". Swap the D/E/F content (and its
#    formatting + row height) between the two rows so FALL BREAK lands on
#    row 14 and the Reshaping class content lands on row 15.
#
# 2) Update performance: strip the stray leading space before each bullet
#    in the "Speed and Performance" class's "Do Before Class" notes.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Stash the current (pre-edit) values ---
$d14 = $ws.Range("D14").Value2
$e14 = $ws.Range("E14").Value2
$f14 = $ws.Range("F14").Value2
$d15 = $ws.Range("D15").Value2

# --- Swap the cell formatting (not just values) between the two rows,
# using copy/PasteSpecial-formats so existing style records are reused
# instead of synthesizing new ones. A scratch cell holds D15's original
# style while D14's moves onto D15.
$ws.Range("D15").Copy()
$ws.Range("Z100").PasteSpecial(-4122)

$ws.Range("D14").Copy()
$ws.Range("D15").PasteSpecial(-4122)

$ws.Range("Z100").Copy()
$ws.Range("D14").PasteSpecial(-4122)

$ws.Range("Z100").ClearFormats()

# E14/F14's formatting moves onto E15/F15 (E15/F15 start out empty/default).
$ws.Range("E14").Copy()
$ws.Range("E15").PasteSpecial(-4122)

$ws.Range("F14").Copy()
$ws.Range("F15").PasteSpecial(-4122)

# --- Now move the values/content to match ---
$ws.Range("C14").ClearContents()

$ws.Range("D14").Value2 = $d15
$ws.Range("E14:F14").Clear()
$ws.Rows.Item(14).AutoFit()

$ws.Range("D15").Value2 = $d14
$ws.Range("E15").Value2 = $e14
$ws.Range("F15").Value2 = $f14
$ws.Rows.Item(15).RowHeight = 99

# --- Update performance: drop the leading space on each bullet line ---
$nl = [char]10
$ws.Range("E19").Value2 = '- `Understanding Performance <performance_understanding.ipynb>`_' + $nl + '- `Improving Performance <performance_solutions.ipynb>`_'

# --- View state: scroll position + active selection ---
$ws.Range("D18").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 5
$win.ScrollColumn = 1
